# Update odds on row 3 and add a new row 6 with a new match,
# as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 values ---
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 5.75
$ws.Range("I3").Value = 12
$ws.Range("L3").Value = 11
$ws.Range("U3").Value = 2.75
$ws.Range("V3").Value = 1.4
$ws.Range("X3").Value = 5
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 6.5
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 10
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 34
$ws.Range("AF3").Value = 126
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 51
$ws.Range("AK3").Value = 151
$ws.Range("AL3").Value = 101
$ws.Range("AM3").Value = 81
$ws.Range("AN3").Value = 3
$ws.Range("AO3").Value = 5.5
$ws.Range("AQ3").Value = 15
$ws.Range("AS3").Value = 201
$ws.Range("AU3").Value = 12
$ws.Range("AV3").Value = 101
$ws.Range("AW3").Value = 12
$ws.Range("AZ3").Value = 351

# --- Add new row 6 with a new match ---
$ws.Range("A6").Value = "AeuJP3G5"
$ws.Range("B6").Value = "24/10/2024"
$ws.Range("C6").Value = "12:15"
$ws.Range("D6").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E6").Value = "Al Orubah"
$ws.Range("F6").Value = "Al Shabab"
$ws.Range("G6").Value = 3.7
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 2.05
$ws.Range("L6").Value = 2.63
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 9
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.75
$ws.Range("S6").Value = 1.44
$ws.Range("T6").Value = 2.63
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.83
$ws.Range("W6").Value = 11
$ws.Range("X6").Value = 19
$ws.Range("Y6").Value = 13
$ws.Range("Z6").Value = 41
$ws.Range("AA6").Value = 34
$ws.Range("AB6").Value = 41
$ws.Range("AC6").Value = 9
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 51
$ws.Range("AG6").Value = 800
$ws.Range("AH6").Value = 7
$ws.Range("AI6").Value = 9
$ws.Range("AJ6").Value = 9
$ws.Range("AK6").Value = 17
$ws.Range("AL6").Value = 17
$ws.Range("AM6").Value = 29
$ws.Range("AN6").Value = 5.5
$ws.Range("AO6").Value = 21
$ws.Range("AP6").Value = 29
$ws.Range("AQ6").Value = 67
$ws.Range("AR6").Value = 101
$ws.Range("AS6").Value = 350
$ws.Range("AT6").Value = 2.63
$ws.Range("AU6").Value = 8.5
$ws.Range("AV6").Value = 51
$ws.Range("AW6").Value = 4
$ws.Range("AX6").Value = 11
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 41
$ws.Range("BA6").Value = 51
$ws.Range("BB6").Value = 151
$ws.Range("BC6").Value = 81
$ws.Range("BD6").Value = 81
